$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Extracted on" timestamp embedded in the study description (A2).
$ws.Range("A2").Value = "This is an interesting study.Extracted on : 2022/09/26 12:26:26"

# Add/apply a left+vertically-centered alignment style to column B's data cells
# (B5:B10), matching the new cellXfs entry introduced for this column.
$r = $ws.Range("B5:B10")
$r.HorizontalAlignment = -4131  # xlLeft
$r.VerticalAlignment = -4108    # xlCenter
